$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (rows 2 through 72) from 45189 to 45190 (date serial +1 day)
$ws.Range("C2:C72").Value = 45190
